# Update sheet title to reflect new "through" date
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Through 2022-07-03"

# Update the row label for July to reflect the new "through" date
$ws.Range("A8").Value = "July (through 07-03)"

# Update the July row (row 8) values for columns C..I (2016..2022)
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 9
$ws.Range("H8").Value = 13
$ws.Range("I8").Value = 23

# Update the Total row (row 9) values for columns C..I (2016..2022)
$ws.Range("C9").Value = 253
$ws.Range("D9").Value = 395
$ws.Range("E9").Value = 362
$ws.Range("F9").Value = 257
$ws.Range("G9").Value = 481
$ws.Range("H9").Value = 773
$ws.Range("I9").Value = 829
